$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is a template: row 1 holds the column headers ("name", "on",
# "axis", "line", "x"). Add two more rows that show the possible values
# for the "on" column ("axis" and "line"), leaving every other column on
# those rows present but blank - matching the other template cells.

# Row 2 -> on: axis
$ws.Cells.Item(2, 1).Value = "'"
$ws.Cells.Item(2, 1).Style = "Normal"
$ws.Cells.Item(2, 2).Value = "axis"
$ws.Cells.Item(2, 3).Value = "'"
$ws.Cells.Item(2, 3).Style = "Normal"
$ws.Cells.Item(2, 4).Value = "'"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "'"
$ws.Cells.Item(2, 5).Style = "Normal"

# Row 3 -> on: line
$ws.Cells.Item(3, 1).Value = "'"
$ws.Cells.Item(3, 1).Style = "Normal"
$ws.Cells.Item(3, 2).Value = "line"
$ws.Cells.Item(3, 3).Value = "'"
$ws.Cells.Item(3, 3).Style = "Normal"
$ws.Cells.Item(3, 4).Value = "'"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "'"
$ws.Cells.Item(3, 5).Style = "Normal"
